# feat: add 2022-Q1 data
#
# Before:  Sheet1 "2021-Q4" | Sheet2 "总计"
# After:   Sheet1 "2021-Q4" | Sheet2 "2022-Q1" (new fund-holding data,
#          replacing the old "总计" sheet content) | Sheet3 "总计" (new
#          summary sheet, gains a 2022-Q1 row ahead of the existing
#          2021-Q4 row)

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. Turn the old "总计" sheet into the new "2022-Q1" fund-holdings sheet
#    and append a brand-new "总计" sheet at the end, mirroring how the
#    sheetIds line up in the target workbook (1, 2, 3 in tab order).
# ---------------------------------------------------------------------
$wsQ4   = $wb.Worksheets.Item(1)
$wsQ1   = $wb.Worksheets.Item(2)
$wsQ1.Name = "2022-Q1"

$wsTotal = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsTotal.Name = "总计"

# ---------------------------------------------------------------------
# 2. Clear whatever the old "总计" sheet had and write the 2022-Q1
#    fund-holding table (header row + 7 data rows, columns A:H).
# ---------------------------------------------------------------------
$wsQ1.Cells.Clear()

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $wsQ1.Cells.Item(1, 2 + $i).Value = $headers[$i]
}

$rows = @(
    @("010054", "万家健康产业混合A", "8.13", "86.63", "3.24", "0.2634", 6),
    @("006132", "万家智造优势混合A", "4.82", "93.70", "2.48", "0.1195", 10),
    @("010055", "万家健康产业混合C", "3.36", "86.63", "3.24", "0.1089", 6),
    @("005108", "圆信永丰双利优选定期开放灵活配置混合", "1.89", "94.60", "3.92", "0.0741", 8),
    @("010434", "红土创新医疗保健股票", "0.75", "92.96", "6.38", "0.0478", 1),
    @("006133", "万家智造优势混合C", "0.52", "93.70", "2.48", "0.0129", 10),
    @("006274", "圆信永丰医药健康混合", "0.18", "93.66", "4.11", "0.0074", 10)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    $excelRow = 2 + $r

    $wsQ1.Cells.Item($excelRow, 1).Value = $r

    # Text-like numeric-looking columns (B..G) must stay text, like the
    # source "2021-Q4" sheet -- force a text number format before typing
    # them in so they aren't auto-coerced to real numbers.
    $textRange = $wsQ1.Range($wsQ1.Cells.Item($excelRow, 2), $wsQ1.Cells.Item($excelRow, 7))
    $textRange.NumberFormat = "@"

    $wsQ1.Cells.Item($excelRow, 2).Value = $row[0]
    $wsQ1.Cells.Item($excelRow, 3).Value = $row[1]
    $wsQ1.Cells.Item($excelRow, 4).Value = $row[2]
    $wsQ1.Cells.Item($excelRow, 5).Value = $row[3]
    $wsQ1.Cells.Item($excelRow, 6).Value = $row[4]
    $wsQ1.Cells.Item($excelRow, 7).Value = $row[5]

    # Rank column (H) is a real number, no special formatting.
    $wsQ1.Cells.Item($excelRow, 8).Value = $row[6]
}

# Match the look of the "2021-Q4" sheet: bold/bordered/centred header row
# and index column -- copy the formatting from the existing sheet so the
# same style is reused instead of creating new ones.
$wsQ4.Range("B1:H1").Copy() | Out-Null
$wsQ1.Range("B1:H1").PasteSpecial($xlPasteFormats) | Out-Null

$wsQ4.Range("A2").Copy() | Out-Null
$wsQ1.Range("A2:A8").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. Populate the brand-new "总计" summary sheet: same layout as before,
#    plus a new 2022-Q1 row ahead of the existing 2021-Q4 row.
# ---------------------------------------------------------------------
$wsTotal.Cells.Item(1, 2).Value = "日期"
$wsTotal.Cells.Item(1, 3).Value = "持有数量(只)"
$wsTotal.Cells.Item(1, 4).Value = "持有市值(亿元)"

$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q1"
$wsTotal.Cells.Item(2, 3).Value = 7
$wsTotal.Cells.Item(2, 4).Value = 0.63

$wsTotal.Cells.Item(3, 1).Value = 1
$wsTotal.Cells.Item(3, 2).Value = "2021-Q4"
$wsTotal.Cells.Item(3, 3).Value = 2
$wsTotal.Cells.Item(3, 4).Value = 0.3

$wsQ4.Range("B1:D1").Copy() | Out-Null
$wsTotal.Range("B1:D1").PasteSpecial($xlPasteFormats) | Out-Null

$wsQ4.Range("A2").Copy() | Out-Null
$wsTotal.Range("A2:A3").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4. Leave the selection on the first sheet, like the source workbook.
# ---------------------------------------------------------------------
$wsQ4.Select() | Out-Null
$wsQ4.Range("A1").Select() | Out-Null
